$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 11 by copying row 6 (its original data moves down to row 11)
$ws.Range("A6:R6").Copy($ws.Range("A11:R11"))

# Row 2: date + volume/price stats updated
$ws.Range("D2").Value = 44396
$ws.Range("J2").Value = 130
$ws.Range("K2").Value = 22000
$ws.Range("M2").Value = 22000
$ws.Range("P2").Value = 1467

# Row 5: date + volume/price stats updated
$ws.Range("D5").Value = 44365
$ws.Range("J5").Value = 580
$ws.Range("K5").Value = 20000
$ws.Range("M5").Value = 21103
$ws.Range("P5").Value = 1407

# Row 6: date + price stats updated (volume unchanged)
$ws.Range("D6").Value = 44398
$ws.Range("K6").Value = 20000
$ws.Range("L6").Value = 20000
$ws.Range("M6").Value = 20000
$ws.Range("P6").Value = 1333

# Row 7: date + volume/price stats updated
$ws.Range("D7").Value = 44391
$ws.Range("J7").Value = 160
$ws.Range("K7").Value = 20000
$ws.Range("L7").Value = 20000
$ws.Range("M7").Value = 20000
$ws.Range("P7").Value = 1333

# Row 8: date + volume/price stats updated
$ws.Range("D8").Value = 44446
$ws.Range("J8").Value = 150
$ws.Range("K8").Value = 22000
$ws.Range("L8").Value = 24000
$ws.Range("M8").Value = 22667
$ws.Range("P8").Value = 1511

# Row 9: date + volume/price stats updated
$ws.Range("D9").Value = 44435
$ws.Range("J9").Value = 140
$ws.Range("K9").Value = 21000
$ws.Range("L9").Value = 23000
$ws.Range("M9").Value = 21714
$ws.Range("P9").Value = 1448

# Row 10: date + volume/price stats updated
$ws.Range("D10").Value = 44399
$ws.Range("J10").Value = 150
$ws.Range("K10").Value = 22000
$ws.Range("L10").Value = 22000
$ws.Range("M10").Value = 22000
$ws.Range("P10").Value = 1467
